# Weekly data refresh: a new observation (most recent week) is inserted at
# row 39, pushing the existing rows 39-130 down to 40-131.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 39 (shifts rows 39..130 -> 40..131,
# and Excel/the engine grows the sheetData dimension to R131 automatically).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = 44544
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 100112029
$ws.Range("G39").Value = "Orégano"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 35
$ws.Range("K39").Value = 8500
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 8729
$ws.Range("N39").Value = "`$/docena de atados"
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 2910
$ws.Range("Q39").Value = 3
$ws.Range("R39").Value = "Hortaliza"
